$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Q3 and R3 to rounded values
$ws.Range("Q3").Value = 818895
$ws.Range("R3").Value = 7382402

# Remove Z3 (Starttid) and AB3 (Sluttid) contents entirely for row 3
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
